# Weekly fruit/vegetable update: insert a new weekly record as row 370,
# pushing all subsequent rows down by one (last row becomes 494).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 370 (shifts 370..493 down to 371..494)
$ws.Rows.Item(370).Insert()

# Populate the newly inserted row 370 with the new weekly record
$ws.Range("A370").Value = 8
$ws.Range("B370").Value = "Terminal La Palmera de La Serena"
$ws.Range("C370").Value = "Coquimbo"
$ws.Range("D370").Value = 45141
$ws.Range("E370").Value = 4
$ws.Range("F370").Value = 100114014
$ws.Range("G370").Value = "Betarraga"
$ws.Range("H370").Value = "Sin especificar"
$ws.Range("I370").Value = "Primera"
$ws.Range("J370").Value = 2000
$ws.Range("K370").Value = 550
$ws.Range("L370").Value = 600
$ws.Range("M370").Value = 575
$ws.Range("N370").Value = "`$/paquete 3 unidades"
$ws.Range("O370").Value = "Provincia del Elquí"
$ws.Range("P370").Value = 192
$ws.Range("Q370").Value = 3
$ws.Range("R370").Value = "Hortaliza"

# Match the date cell formatting used by the rest of column D
$ws.Range("D370").NumberFormat = $ws.Range("D371").NumberFormat
